$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "average" label + AVERAGE formulas for the two data blocks (C3:C22 and L3:L22)
$ws.Range("B23").Value = "average"
$ws.Range("C23").Formula = "=AVERAGE(C3:C22)"

$ws.Range("K23").Value = "average"
$ws.Range("L23").Formula = "=AVERAGE(L3:L22)"

# Match the saved selection state from the diff
$ws.Range("L23").Select()
